# Rename all 50 worksheets (in their current tab order) to the new
# "summXXXXXXXX" names, per the commit "Add importance ranking and shap
# by postcode for distance scripts".

$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ18051333",
    "summ18301228",
    "summ18646257",
    "summ18961913",
    "summ19321655",
    "summ19669670",
    "summ20060197",
    "summ20420754",
    "summ20836391",
    "summ21220668",
    "summ21549454",
    "summ21893488",
    "summ22243963",
    "summ22556568",
    "summ22893590",
    "summ23230309",
    "summ23553349",
    "summ23909955",
    "summ24334345",
    "summ24676259",
    "summ25011949",
    "summ25356064",
    "summ25680264",
    "summ26067927",
    "summ26486300",
    "summ26820752",
    "summ27179972",
    "summ27598243",
    "summ27981245",
    "summ28389382",
    "summ28747031",
    "summ29141659",
    "summ29478752",
    "summ29910629",
    "summ30254390",
    "summ30578083",
    "summ30925147",
    "summ31269475",
    "summ31626417",
    "summ31999249",
    "summ32389179",
    "summ32850262",
    "summ33219269",
    "summ33630548",
    "summ34018298",
    "summ34375043",
    "summ34715115",
    "summ35032202",
    "summ35365126",
    "summ35686838"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $wb.Worksheets.Item($i).Name = $newNames[$i - 1]
}
